$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Regionalverband Saarbrücken / Dr. Jäger Privatkliniken)
$ws.Range("E4").Value = 49.2371047
$ws.Range("F4").Value = 6.9925617
$ws.Range("G4").Value = 8.66

# Row 9 (Merzig-Wadern / Klinikum Merzig gGmbH)
$ws.Range("E9").Value = 49.4571859
$ws.Range("F9").Value = 6.6315777
$ws.Range("G9").Value = 11.96
$ws.Range("H9").Value = 11.2

# Row 10 (Merzig-Wadern / Klinikum Merzig gGmbH)
$ws.Range("E10").Value = 49.4571859
$ws.Range("F10").Value = 6.6315777
$ws.Range("G10").Value = 24.09
$ws.Range("H10").Value = 15.07
